$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stop names (column B, rows 2-6)
$ws.Range("B2").Value = "Alipurduar"
$ws.Range("B3").Value = "Alipurduar Chowpathy"
$ws.Range("B4").Value = "Alipurduar Birpara"
$ws.Range("B5").Value = "Sonapur"
$ws.Range("B6").Value = "Falakata"

# Update latitude (col C), longitude (col D) and distance_km (col E) values
$ws.Range("C2").Value = 26.4799
$ws.Range("D2").Value = 89.5355

$ws.Range("C3").Value = 26.48083
$ws.Range("D3").Value = 89.526
$ws.Range("E3").Value = 1.5

$ws.Range("C4").Value = 26.48281
$ws.Range("D4").Value = 89.50897
$ws.Range("E4").Value = 2.2

$ws.Range("C5").Value = 26.494
$ws.Range("D5").Value = 89.368
$ws.Range("E5").Value = 10.5

$ws.Range("C6").Value = 26.5193
$ws.Range("D6").Value = 89.202
$ws.Range("E6").Value = 15

# Widen columns B, C, D to fit the longer stop names now in use
$ws.Columns.Item(2).ColumnWidth = 22.833333333333336
$ws.Columns.Item(3).ColumnWidth = 24.666666666666668
$ws.Columns.Item(4).ColumnWidth = 32.33333333333333

# Add a new blank formatted row below the data (mirrors formatting of row above)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Move the active selection to C3
$ws.Range("C3").Select()
